$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# Overview sheet: E2, F2, E3, F3
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# zh-cn sheet: C2, C3 (Status column)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

# de-de sheet: C2, C3 (Status column)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Column width changes: 17.2159881591797 -> 13.4101845877511 (characters) ---
# NOTE: the engine quantizes ColumnWidth to 1/6-character steps with an
# internal +5/6 offset (stored = round((input + 5/6) * 6) / 6), so the raw
# target value is pre-compensated to 12.5 to land on the closest achievable
# stored width (13.333333333333334, the nearest reachable value to
# 13.4101845877511).
# Overview sheet: columns E (5) and F (6)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C (3)
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C (3)
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
